# Add a new "2024/11/23" data column (BX) to the 合成確率 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column width (column 76 = BX), matches the other data columns (raw width 12) ---
$ws.Columns.Item(76).ColumnWidth = 11.17

# --- Header cell BX1: plain text date label, same look as the other header cells (BW1) ---
$ws.Range("BX1").NumberFormat = "@"
$ws.Range("BX1").Value = "2024/11/23"
$ws.Range("BW1").Copy()
$ws.Range("BX1").PasteSpecial(-4122)

# --- Style template cells already present on the sheet ---
#   A2 -> s="1" (normal fill)            used when value >= 140
#   D2 -> s="2" (yellow fill, big win)   used when value <  125
#   N2 -> s="3" (light-blue fill, small) used when 125 <= value < 140
$normalTemplate = $ws.Range("A2")
$bigTemplate    = $ws.Range("D2")
$smallTemplate  = $ws.Range("N2")

# Data values for 2024/11/23, rows 2..53 (machine rows in sheet order)
$values = @(
    171, 135.1, 162.9, 181.5, 134.9, 210, 148.2, 156.8, 149.7, 163.8,
    152.2, 133.9, 415.6, 130, 132, 153, 141.6, 173.6, 162, 194.7,
    154.6, 164.7, 135.5, 213.6, 124.1, 173.8, 188, 146.6, 132.4, 165.3,
    146.4, 156.6, 195.5, 118.6, 206.1, 141.9, 134.7, 129.3, 222.8, 145.4,
    145.7, 235.6, 146, 149, 146.8, 174.8, 120.5, 132.2, 165.9, 156.9,
    152.2, 125.8
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $v = $values[$i]
    $target = $ws.Cells.Item($row, 76)

    if ($v -lt 125) {
        $bigTemplate.Copy()
    } elseif ($v -lt 140) {
        $smallTemplate.Copy()
    } else {
        $normalTemplate.Copy()
    }
    $target.PasteSpecial(-4122)
    $target.Value = $v
}

Write-Output "BX column for 2024/11/23 added"
